# Applies the documentation update described by the commit:
# "added some more documentation to the project"
#
# 1) Extends the paragraph ending in "lehetnek." with a new sentence about
#    the AnimationTimer handler function.
# 2) Splits the run "(Game) fuggvenyt, ami olyan mint az osszes " into
#    "(" + "Game) fuggvenyt, ami olyan mint az osszes " (no text change).
# 3) Fixes capitalization "(Ami" -> "(ami" inside the " attributumat (Ami egy "
#    run.
# 4) Adds a brand-new paragraph at the end describing how invaders shoot.

$d = $word.ActiveDocument

# Helper: insert `$text` at `$startPos` (a numeric character offset) using
# the character formatting captured in `$ft` (a Word `FormattedText` Range
# obtained from a correctly-formatted neighbouring run). Returns the
# position immediately after the inserted text so callers can chain calls.
function Insert-FormattedRun($startPos, $text, $ft) {
    $beforeEnd = $d.Content.End
    $ip = $d.Range($startPos, $startPos)
    $ip.FormattedText = $ft
    $afterEnd = $d.Content.End
    $insertedLen = $afterEnd - $beforeEnd
    $newRange = $d.Range($startPos, $startPos + $insertedLen)
    $newRange.Text = $text
    return $startPos + $text.Length
}

# ---------------------------------------------------------------------
# Change 1: " Ezt egy AnimationTimer handler fuggvenyeben valositottam
# meg, hogy a jatek folyamatosan fusson." after "lehetnek."
# ---------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Execute("lehetnek.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$ft1 = $find1.FormattedText
$pos = $find1.End

$pos = Insert-FormattedRun $pos " Ezt egy " $ft1
$pos = Insert-FormattedRun $pos "AnimationTimer" $ft1
$pos = Insert-FormattedRun $pos " " $ft1
$pos = Insert-FormattedRun $pos "handler" $ft1
$pos = Insert-FormattedRun $pos " függvényében valósítottam meg, hogy a játék folyamatosan fusson." $ft1

# ---------------------------------------------------------------------
# Change 2: split "(Game) fuggvenyt, ami olyan mint az osszes " into two
# runs "(" / "Game) ...". Text itself is unchanged, so a plain
# Find/Replace (same text back) is sufficient to normalise on content;
# the surrounding formatting is preserved automatically since it is an
# in-place replace, not an insert.
# ---------------------------------------------------------------------
$find2 = $d.Content
$find2.Find.Execute("(Game) függvényt, ami olyan mint az összes ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$find2.Text = "(Game) függvényt, ami olyan mint az összes "

# ---------------------------------------------------------------------
# Change 3: " attributumat (Ami egy " -> " attributumat (ami egy "
# (lower-cases the "Ami", also later split into multiple runs upstream -
# the visible text result is what matters here).
# ---------------------------------------------------------------------
$find3 = $d.Content
$find3.Find.Execute(" attribútumát (Ami egy ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$find3.Text = " attribútumát (ami egy "

# ---------------------------------------------------------------------
# Change 4: brand-new paragraph at the end of the document describing
# how invaders fire.
# ---------------------------------------------------------------------
$find4 = $d.Content
$find4.Find.Execute("és stand értékeit veheti fel).", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$find4.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$pRange = $newPara.Range
$pRange.Collapse(1)
$pRange.InsertAfter("Az invaderek lövéséhez az AnimationTimer után létrehoztam egy Timeline-t, aminek feladata pár másodpercenként végigmenni az invadereken, amiken kis valószínűséggel meghívni a shoot() függvényüket.")
